# Add new archetype attributes relating home distance to WoS (Workplace/School)
# and POI (Point of Interest): dist_wos_mu, dist_wos_sigma, dist_poi_mu, dist_poi_sigma

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New headers in row 1, columns AC:AF
$ws.Range("AC1").Value = "dist_wos_mu"
$ws.Range("AD1").Value = "dist_wos_sigma"
$ws.Range("AE1").Value = "dist_poi_mu"
$ws.Range("AF1").Value = "dist_poi_sigma"

# New data values for the 5 archetype rows (2-6): dist_wos_mu, dist_wos_sigma, dist_poi_mu, dist_poi_sigma
$ws.Cells.Item(2, 29).Value = 14000
$ws.Cells.Item(2, 30).Value = 500
$ws.Cells.Item(2, 31).Value = 2000
$ws.Cells.Item(2, 32).Value = 500

$ws.Cells.Item(3, 29).Value = 14000
$ws.Cells.Item(3, 30).Value = 500
$ws.Cells.Item(3, 31).Value = 2000
$ws.Cells.Item(3, 32).Value = 500

$ws.Cells.Item(4, 29).Value = 2000
$ws.Cells.Item(4, 30).Value = 500
$ws.Cells.Item(4, 31).Value = 2000
$ws.Cells.Item(4, 32).Value = 500

$ws.Cells.Item(5, 29).Value = 2000
$ws.Cells.Item(5, 30).Value = 500
$ws.Cells.Item(5, 31).Value = 2000
$ws.Cells.Item(5, 32).Value = 500

$ws.Cells.Item(6, 29).Value = 14000
$ws.Cells.Item(6, 30).Value = 500
$ws.Cells.Item(6, 31).Value = 2000
$ws.Cells.Item(6, 32).Value = 500

# Size the new columns to their best-fit widths (as Excel would compute for this content)
$ws.Range("AC1").ColumnWidth = 8405.0 / 768.0
$ws.Range("AD1").ColumnWidth = 79.0 / 6.0
$ws.Range("AE1").ColumnWidth = 1973.0 / 192.0
$ws.Range("AF1").ColumnWidth = 9599.0 / 768.0

# Update selection to match final state
$ws.Range("AB22").Select() | Out-Null
